# Update the cryptos price/volume table with the latest scraped values.
# Column D ("Price") values are forced to literal text via a leading
# apostrophe so Excel does not auto-convert look-alike numbers (e.g.
# "597.43", "6.40", "1.00") into numeric cells - they must stay plain text,
# matching the original inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.560.74"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "'3.020.62"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'597.43"
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("D6").Value = "'149.83"
$ws.Range("E6").Value = "  +5.95%  "
$ws.Range("D8").Value = "'3.018.87"
$ws.Range("E8").Value = "  +1.56%  "
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").Value = "'6.40"
$ws.Range("E10").Value = "  +11.37%  "
$ws.Range("D11").Value = "'0.150"
$ws.Range("E11").Value = "  +4.17%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("E13").Value = "  +3.56%  "
$ws.Range("D14").Value = "'34.53"
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("E15").Value = "  +2.63%  "
$ws.Range("D16").Value = "'3.521.68"
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("D17").Value = "'62.548.48"
$ws.Range("E17").Value = "  +1.75%  "
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "'3.025.82"
$ws.Range("E19").Value = "  +2.43%  "
$ws.Range("D20").Value = "'448.06"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "'14.17"
$ws.Range("E21").Value = "  +2.50%  "
$ws.Range("D22").Value = "'0.689"
$ws.Range("E22").Value = "  +1.27%  "
$ws.Range("D23").Value = "'7.44"
$ws.Range("E23").Value = "  +2.02%  "
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("D25").Value = "'10.85"
$ws.Range("E25").Value = "  +12.76%  "
$ws.Range("E26").Value = "  +4.80%  "
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  +3.88%  "
$ws.Range("D30").Value = "'7.28"
$ws.Range("E30").Value = "  +6.10%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  +4.35%  "
$ws.Range("D33").Value = "'27.53"
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("E34").Value = "  +3.04%  "
$ws.Range("D35").Value = "'0.0₃0852"
$ws.Range("E35").Value = "  +10.63%  "
$ws.Range("E36").Value = "  +1.84%  "
$ws.Range("D37").Value = "'5.84"
$ws.Range("E37").Value = "  +3.23%  "
$ws.Range("D38").Value = "'3.04"
$ws.Range("E38").Value = "  +9.23%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "'50.06"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("D42").Value = "'0.124"
$ws.Range("E42").Value = "  +4.37%  "
$ws.Range("D43").Value = "'0.286"
$ws.Range("E43").Value = "  +9.48%  "
$ws.Range("D44").Value = "'391.22"
$ws.Range("E44").Value = "  +1.19%  "
$ws.Range("D45").Value = "'40.14"
$ws.Range("E45").Value = "  +8.55%  "
$ws.Range("D46").Value = "'0.0354"
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("D47").Value = "'2.738.26"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("D48").Value = "'132.91"
$ws.Range("E48").Value = "  +2.45%  "
$ws.Range("E50").Value = "  +1.77%  "
$ws.Range("E51").Value = "  +0.15%  "
